$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the "Last status check" timestamp text in F1
$ws.Range("F1").Value = "Last status check on: 01.02.2022 06:15"

# 2. D4 used to be inline text "+0.4"; it becomes a real numeric value 0.4
$ws.Range("D4").Value = 0.4

# 3. E4 used to be inline text "2022-02-01 06:00:13"; it becomes a real
#    date/time value (Excel serial number 44593.25015046296), formatted
#    with the same date/time style already used by sibling cells (E2, E3, ...)
$ws.Range("E4").Value = 44593.25015046296
$ws.Range("E4").NumberFormat = "YYYY-MM-DD HH:MM:SS"
